# Edit script: add more diseases to SymptomTable
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing "Suara Serak" label (renamed to lower-case "serak") that is
# referenced from several existing rows (B7, B52, B86, B117).
$ws.Cells.Item(7, 2).Value = "Suara serak"
$ws.Cells.Item(52, 2).Value = "Suara serak"
$ws.Cells.Item(86, 2).Value = "Suara serak"
$ws.Cells.Item(117, 2).Value = "Suara serak"

# New rows (134-233) describing additional diseases/symptoms.
$newRows = @(
    @{Row=134; A='Dyspepsia'; B='Nyeri ulu hati'; C=''; D='Sangat sering'},
    @{Row=135; A='Dyspepsia'; B='Nyeri memberat setelah makan'; C=''; D='Jarang'},
    @{Row=136; A='Dyspepsia'; B='Regurgitasi (Rasa pahit/asam di mulut)'; C=''; D='Jarang'},
    @{Row=137; A='Dyspepsia'; B='Perut kembung'; C=''; D='Sering'},
    @{Row=138; A='Dyspepsia'; B='Mual'; C=''; D='Sering'},
    @{Row=139; A='Dyspepsia'; B='Anorexia (hilang nafsu makan)'; C=''; D='Kadang'},
    @{Row=140; A='GERD'; B='Nyeri ulu hati'; C=''; D='Sangat sering'},
    @{Row=141; A='GERD'; B='Nyeri memberat setelah makan'; C=''; D='Sangat sering'},
    @{Row=142; A='GERD'; B='Rasa panas di dada'; C=''; D='Sering'},
    @{Row=143; A='GERD'; B='Regurgitasi (Rasa pahit/asam di mulut)'; C=''; D='Sering'},
    @{Row=144; A='GERD'; B='Perut kembung'; C=''; D='Jarang'},
    @{Row=145; A='GERD'; B='Mual'; C=''; D='Sering'},
    @{Row=146; A='GERD'; B='Muntah'; C=''; D='Jarang'},
    @{Row=147; A='GERD'; B='Batuk'; C=''; D='Kadang'},
    @{Row=148; A='GERD'; B='Suara serak'; C=''; D='Kadang'},
    @{Row=149; A='GERD'; B='Sesak napas'; C=''; D='Kadang'},
    @{Row=150; A='GERD'; B='Anorexia (hilang nafsu makan)'; C=''; D='Kadang'},
    @{Row=151; A='Appendicitis Akut'; B='Nyeri perut kanan bawah'; C=''; D='Sering'},
    @{Row=152; A='Appendicitis Akut'; B='Nyeri ulu hati'; C=''; D='Kadang'},
    @{Row=153; A='Appendicitis Akut'; B='Nyeri perut kanan atas'; C=''; D='Jarang'},
    @{Row=154; A='Appendicitis Akut'; B='Perut kembung'; C=''; D='Sering'},
    @{Row=155; A='Appendicitis Akut'; B='Mual'; C=''; D='Sering'},
    @{Row=156; A='Appendicitis Akut'; B='Muntah'; C=''; D='Kadang'},
    @{Row=157; A='Appendicitis Akut'; B='Demam'; C=''; D='Sering'},
    @{Row=158; A='Appendicitis Akut'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=159; A='Appendicitis Akut'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=160; A='Appendicitis Akut'; B='Diare (BAB encer ≥3x/hari)'; C=''; D='Jarang'},
    @{Row=161; A='Appendicitis Akut'; B='Konstipasi'; C=''; D='Jarang'},
    @{Row=162; A='Gastroenteritis Akut'; B='Nyeri perut (tidak terlokalisir)'; C=''; D='Sering'},
    @{Row=163; A='Gastroenteritis Akut'; B='Perut kembung'; C=''; D='Kadang'},
    @{Row=164; A='Gastroenteritis Akut'; B='Mual'; C=''; D='Sering'},
    @{Row=165; A='Gastroenteritis Akut'; B='Muntah'; C=''; D='Kadang'},
    @{Row=166; A='Gastroenteritis Akut'; B='Demam'; C=''; D='Sering'},
    @{Row=167; A='Gastroenteritis Akut'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=168; A='Gastroenteritis Akut'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=169; A='Gastroenteritis Akut'; B='Diare (BAB encer ≥3x/hari)'; C=''; D='Sering'},
    @{Row=170; A='Gastroenteritis Akut'; B='Konstipasi'; C=''; D='Jarang'},
    @{Row=171; A='Gastroenteritis Akut'; B='Myalgia (nyeri tubuh)'; C=''; D='Kadang'},
    @{Row=172; A='Demam Tifoid'; B='Nyeri perut kanan bawah'; C=''; D='Jarang'},
    @{Row=173; A='Demam Tifoid'; B='Nyeri ulu hati'; C=''; D='Jarang'},
    @{Row=174; A='Demam Tifoid'; B='Nyeri perut (tidak terlokalisir)'; C=''; D='Jarang'},
    @{Row=175; A='Demam Tifoid'; B='Perut kembung'; C=''; D='Kadang'},
    @{Row=176; A='Demam Tifoid'; B='Mual'; C=''; D='Kadang'},
    @{Row=177; A='Demam Tifoid'; B='Muntah'; C=''; D='Kadang'},
    @{Row=178; A='Demam Tifoid'; B='Demam'; C=''; D='Sangat sering'},
    @{Row=179; A='Demam Tifoid'; B='Durasi demam'; C='> 7 hari'; D=''},
    @{Row=180; A='Demam Tifoid'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=181; A='Demam Tifoid'; B='Diare (BAB encer ≥3x/hari)'; C=''; D='Kadang'},
    @{Row=182; A='Demam Tifoid'; B='Konstipasi'; C=''; D='Kadang'},
    @{Row=183; A='Demam Tifoid'; B='Nyeri kepala'; C=''; D='Kadang'},
    @{Row=184; A='Demam Tifoid'; B='Myalgia (nyeri tubuh)'; C=''; D='Kadang'},
    @{Row=185; A='Pankreatitis Akut'; B='Nyeri ulu hati'; C=''; D='Sangat sering'},
    @{Row=186; A='Pankreatitis Akut'; B='Nyeri perut (tidak terlokalisir)'; C=''; D='Jarang'},
    @{Row=187; A='Pankreatitis Akut'; B='Nyeri memberat setelah makan'; C=''; D='Sangat sering'},
    @{Row=188; A='Pankreatitis Akut'; B='Perut kembung'; C=''; D='Kadang'},
    @{Row=189; A='Pankreatitis Akut'; B='Mual'; C=''; D='Sering'},
    @{Row=190; A='Pankreatitis Akut'; B='Muntah'; C=''; D='Sering'},
    @{Row=191; A='Pankreatitis Akut'; B='Sesak napas'; C=''; D='Jarang'},
    @{Row=192; A='Pankreatitis Akut'; B='Demam'; C=''; D='Kadang'},
    @{Row=193; A='Pankreatitis Akut'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=194; A='Pankreatitis Akut'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=195; A='Pankreatitis Akut'; B='Diare (BAB encer ≥3x/hari)'; C=''; D='Jarang'},
    @{Row=196; A='Pankreatitis Akut'; B='Konstipasi'; C=''; D='jarang'},
    @{Row=197; A='Cholecystitis'; B='Nyeri perut (tidak terlokalisir)'; C=''; D='Jarang'},
    @{Row=198; A='Cholecystitis'; B='Nyeri perut kanan atas'; C=''; D='Sangat sering'},
    @{Row=199; A='Cholecystitis'; B='Nyeri memberat setelah makan'; C=''; D='Sering'},
    @{Row=200; A='Cholecystitis'; B='Perut kembung'; C=''; D='Jarang'},
    @{Row=201; A='Cholecystitis'; B='Mual'; C=''; D='Kadang'},
    @{Row=202; A='Cholecystitis'; B='Muntah'; C=''; D='Kadang'},
    @{Row=203; A='Cholecystitis'; B='Demam'; C=''; D='Kadang'},
    @{Row=204; A='Cholecystitis'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=205; A='Cholecystitis'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=206; A='Cholecystitis'; B='Kekuningan pada Tubuh (mata/badan)'; C=''; D='Sering'},
    @{Row=207; A='Hepatitis'; B='Nyeri perut kanan atas'; C=''; D='Kadang'},
    @{Row=208; A='Hepatitis'; B='Nyeri memberat setelah makan'; C=''; D='Kadang'},
    @{Row=209; A='Hepatitis'; B='Perut kembung'; C=''; D='Sering'},
    @{Row=210; A='Hepatitis'; B='Perut membuncit'; C=''; D='Kadang'},
    @{Row=211; A='Hepatitis'; B='Mual'; C=''; D='Kadang'},
    @{Row=212; A='Hepatitis'; B='Muntah'; C=''; D='Kadang'},
    @{Row=213; A='Hepatitis'; B='Sesak napas'; C=''; D='Jarang'},
    @{Row=214; A='Hepatitis'; B='Demam'; C=''; D='Kadang'},
    @{Row=215; A='Hepatitis'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=216; A='Hepatitis'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sering'},
    @{Row=217; A='Hepatitis'; B='Penurunan berat badan'; C=''; D='Kadang'},
    @{Row=218; A='Hepatitis'; B='Myalgia (nyeri tubuh)'; C=''; D='Kadang'},
    @{Row=219; A='Hepatitis'; B='Kekuningan pada Tubuh (mata/badan)'; C=''; D='Sangat sering'},
    @{Row=220; A='Hepatitis'; B='Pucat (Anemia)'; C=''; D='jarang'},
    @{Row=221; A='Ascariasis'; B='Nyeri perut (tidak terlokalisir)'; C=''; D='Jarang'},
    @{Row=222; A='Ascariasis'; B='Perut kembung'; C=''; D='Kadang'},
    @{Row=223; A='Ascariasis'; B='Perut membuncit'; C=''; D='Kadang'},
    @{Row=224; A='Ascariasis'; B='Mual'; C=''; D='Sering'},
    @{Row=225; A='Ascariasis'; B='Muntah'; C=''; D='Kadang'},
    @{Row=226; A='Ascariasis'; B='Batuk'; C=''; D='Jarang'},
    @{Row=227; A='Ascariasis'; B='Demam'; C=''; D='Jarang'},
    @{Row=228; A='Ascariasis'; B='Durasi demam'; C='< 7 hari'; D=''},
    @{Row=229; A='Ascariasis'; B='Anorexia (hilang nafsu makan)'; C=''; D='Sangat sering'},
    @{Row=230; A='Ascariasis'; B='Penurunan berat badan'; C=''; D='sangat sering'},
    @{Row=231; A='Ascariasis'; B='Diare (BAB encer ≥3x/hari)'; C=''; D='Jarang'},
    @{Row=232; A='Ascariasis'; B='Konstipasi'; C=''; D='Kadang'},
    @{Row=233; A='Ascariasis'; B='Pucat (Anemia)'; C=''; D='Sangat sering'}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    if ($r.C -ne '') {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    if ($r.D -ne '') {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
}

# Update the selection to match the new bottom of the sheet (the frozen
# header pane at row 1 is left untouched).
$ws.Range("A233").Select()
